$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (August): H9 157 -> 158
$ws.Range("H9").Value = 158

# Row 10 (September): update label and counts
$ws.Range("A10").Value = "September (through 09-04)"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 11
$ws.Range("G10").Value = 11
$ws.Range("H10").Value = 17

# Row 11 (Total): update sums
$ws.Range("B11").Value = 196
$ws.Range("C11").Value = 386
$ws.Range("D11").Value = 563
$ws.Range("E11").Value = 496
$ws.Range("F11").Value = 360
$ws.Range("G11").Value = 795
$ws.Range("H11").Value = 1088
